# "push changes to fix gender"
#
# The FHIR IG publisher re-ran and regenerated this ValueSet export a bit
# later than before (new run timestamp), and this time it did not emit the
# separate "Include ValueSet #0" sheet (the explicit ValueSet-URL include).
# Only the code-system include ("Include #1") remains, and it is renumbered
# to "Include #0" since it is now the only include sheet.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation timestamp (Date row) -------------
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Cells.Item(8, 2).Value = "2025-06-23T13:45:54+01:00"

# --- Drop the "Include ValueSet #0" sheet (ValueSet URL include) ----------
$includeValueSetSheet = $wb.Worksheets.Item("Include ValueSet #0")
$null = $includeValueSetSheet.Delete()

# --- The remaining include sheet becomes "Include #0" ----------------------
$includeSheet = $wb.Worksheets.Item("Include #1")
$includeSheet.Name = "Include #0"
